$wb = $excel.ActiveWorkbook

# Update the A2 "co2" price value on each year sheet with the latest server results.
# Sheet "2035" is left untouched (value unchanged in source data).
$wb.Worksheets.Item("2025").Range("A2").Value = 57
$wb.Worksheets.Item("2030").Range("A2").Value = 195
$wb.Worksheets.Item("2040").Range("A2").Value = 355
$wb.Worksheets.Item("2045").Range("A2").Value = 355
$wb.Worksheets.Item("2050").Range("A2").Value = 355
